# Fix naive component forecaster bug - update forecast error values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

$data = @{
    2  = @(-0.1829735462431551, 0.6694626356775797, 0.8372369824267251, 0.9150065477507389, 0.905446197782048, 51)
    3  = @(0.1263938744343593, 0.8614631615508216, 1.473077372104747, 1.213703988666408, 1.219359969539984, 50)
    4  = @(-0.09485759795785621, 0.6206824703513352, 0.6573762981246858, 0.8107874555792571, 0.813563898042168, 49)
    5  = @(0.05113506990239874, 0.7519750174941408, 1.034754511348446, 1.017228839223725, 1.026693788557885, 48)
    6  = @(-0.03772088114894339, 0.5996523225882197, 0.6704444970345537, 0.8188067519473405, 0.8267802493317774, 47)
    7  = @(0.1077066685141252, 0.790586417127402, 1.154731043345991, 1.074584125764936, 1.080987147395736, 46)
    8  = @(-0.03868616839471213, 0.5498893294951147, 0.5944205642755658, 0.7709867471465159, 0.7787165651551944, 45)
    9  = @(0.1749977985370214, 0.7501568359084184, 1.058248496194531, 1.028712056988996, 1.025437714082251, 44)
    10 = @(0.02191909566903653, 0.5542941723622423, 0.633576040990445, 0.7959748997238826, 0.8050896216323691, 43)
    11 = @(0.2576676093386086, 0.7888566722486881, 1.23948852930275, 1.113323191756441, 1.096224442490078, 42)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
